$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F ("From RCSB") duplicating the "Is model" column (E) values
$ws.Range("F1").Value = "From RCSB"
$ws.Range("F2").Value = $ws.Range("E2").Value2
$ws.Range("F3").Value = $ws.Range("E3").Value2
$ws.Range("F4").Value = $ws.Range("E4").Value2
$ws.Range("F5").Value = $ws.Range("E5").Value2

# Update the sheet view's selection to span the new used range
[void]$ws.Range("A1:F5").Select()

# Match the recorded page setup tweak (portrait orientation)
$ws.PageSetup.Orientation = 1
